$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1901.5
$ws.Range("I28").Value = 2267.6667
$ws.Range("J28").Value = 803
$ws.Range("K28").Value = 2267.6667
$ws.Range("L28").Value = 803
$ws.Range("M28").Value = -1782.6667
$ws.Range("N28").Value = -1773
$ws.Range("H111").Value = 9980.684999999999
$ws.Range("I111").Value = 2284.3333
$ws.Range("J111").Value = 16907.4
$ws.Range("K111").Value = 6852.999899999999
$ws.Range("L111").Value = 50722.2
$ws.Range("M111").Value = -3785.999899999999
$ws.Range("N111").Value = -56856.2
$ws.Range("H112").Value = 4419.615
$ws.Range("J112").Value = 4419.615
$ws.Range("L112").Value = 13258.845
$ws.Range("N112").Value = -15474.845
$ws.Range("H113").Value = 3431.4285
$ws.Range("I113").Value = 2505
$ws.Range("J113").Value = 3802
$ws.Range("K113").Value = 2505
$ws.Range("L113").Value = 3802
$ws.Range("M113").Value = 749
$ws.Range("N113").Value = -10310
$ws.Range("H125").Value = 6171.143
$ws.Range("I125").Value = 2666
$ws.Range("J125").Value = 7573.2
$ws.Range("K125").Value = 23994
$ws.Range("L125").Value = 68158.8
$ws.Range("M125").Value = -21534
$ws.Range("N125").Value = -73078.8
$ws.Range("H138").Value = 3459.0137
$ws.Range("I138").Value = 1890.0333
$ws.Range("J138").Value = 4553.6514
$ws.Range("K138").Value = 5670.0999
$ws.Range("L138").Value = 13660.9542
$ws.Range("M138").Value = -530.0999000000002
$ws.Range("N138").Value = -23940.9542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2535.0789
$ws.Range("I61").Value = 2388.1936
$ws.Range("J61").Value = 3185.5715
$ws.Range("K61").Value = 2388.1936
$ws.Range("L61").Value = 3185.5715
$ws.Range("M61").Value = -2176.1936
$ws.Range("N61").Value = -3609.5715
$ws.Range("H88").Value = 2718.889
$ws.Range("I88").Value = 1996.6666
$ws.Range("J88").Value = 3080
$ws.Range("K88").Value = 1996.6666
$ws.Range("L88").Value = 3080
$ws.Range("M88").Value = -1590.6666
$ws.Range("N88").Value = -3892
$ws.Range("H91").Value = 2718.889
$ws.Range("I91").Value = 1996.6666
$ws.Range("J91").Value = 3080
$ws.Range("K91").Value = 1996.6666
$ws.Range("L91").Value = 3080
$ws.Range("M91").Value = -592.6666
$ws.Range("N91").Value = -5888
$ws.Range("H122").Value = 5157.613
$ws.Range("I122").Value = 5751.68
$ws.Range("J122").Value = 2682.3333
$ws.Range("K122").Value = 17255.04
$ws.Range("L122").Value = 8046.999899999999
$ws.Range("M122").Value = -14805.04
$ws.Range("N122").Value = -12946.9999
$ws.Range("H136").Value = 2535.0789
$ws.Range("I136").Value = 2388.1936
$ws.Range("J136").Value = 3185.5715
$ws.Range("K136").Value = 7164.5808
$ws.Range("L136").Value = 9556.7145
$ws.Range("M136").Value = -4614.5808
$ws.Range("N136").Value = -14656.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2469938.8
$ws.Range("I80").Value = 6174689.5
$ws.Range("J80").Value = 105.22222
$ws.Range("K80").Value = 6174689.5
$ws.Range("L80").Value = 105.22222
$ws.Range("M80").Value = -6173691.5
$ws.Range("N80").Value = -2101.22222
$ws.Range("H83").Value = 2469938.8
$ws.Range("I83").Value = 6174689.5
$ws.Range("J83").Value = 105.22222
$ws.Range("K83").Value = 30873447.5
$ws.Range("L83").Value = 526.1111
$ws.Range("M83").Value = -30868455.5
$ws.Range("N83").Value = -10510.1111
$ws.Range("H86").Value = 168847.75
$ws.Range("I86").Value = 2398.25
$ws.Range("J86").Value = 501746.75
$ws.Range("K86").Value = 2398.25
$ws.Range("L86").Value = 501746.75
$ws.Range("M86").Value = -1275.25
$ws.Range("N86").Value = -503992.75
$ws.Range("H89").Value = 168847.75
$ws.Range("I89").Value = 2398.25
$ws.Range("J89").Value = 501746.75
$ws.Range("K89").Value = 11991.25
$ws.Range("L89").Value = 2508733.75
$ws.Range("M89").Value = -6375.25
$ws.Range("N89").Value = -2519965.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2123.152
$ws.Range("I31").Value = 1211.2916
$ws.Range("J31").Value = 3117.9092
$ws.Range("K31").Value = 1211.2916
$ws.Range("L31").Value = 3117.9092
$ws.Range("M31").Value = -916.2916
$ws.Range("N31").Value = -3707.9092
$ws.Range("H34").Value = 2123.152
$ws.Range("I34").Value = 1211.2916
$ws.Range("J34").Value = 3117.9092
$ws.Range("K34").Value = 1211.2916
$ws.Range("L34").Value = 3117.9092
$ws.Range("M34").Value = -1009.2916
$ws.Range("N34").Value = -3521.9092
$ws.Range("H86").Value = 10339.714
$ws.Range("I86").Value = 11029.667
$ws.Range("J86").Value = 6200
$ws.Range("K86").Value = 11029.667
$ws.Range("L86").Value = 6200
$ws.Range("M86").Value = -9906.666999999999
$ws.Range("N86").Value = -8446
$ws.Range("H89").Value = 10339.714
$ws.Range("I89").Value = 11029.667
$ws.Range("J89").Value = 6200
$ws.Range("K89").Value = 55148.335
$ws.Range("L89").Value = 31000
$ws.Range("M89").Value = -49532.335
$ws.Range("N89").Value = -42232
$ws.Range("H122").Value = 6728.25
$ws.Range("I122").Value = 9322.4
$ws.Range("J122").Value = 2404.6667
$ws.Range("K122").Value = 27967.2
$ws.Range("L122").Value = 7214.000100000001
$ws.Range("M122").Value = -25517.2
$ws.Range("N122").Value = -12114.0001
$ws.Range("H132").Value = 3322.1333
$ws.Range("I132").Value = 2761.0833
$ws.Range("J132").Value = 5566.3335
$ws.Range("K132").Value = 8283.249899999999
$ws.Range("L132").Value = 16699.0005
$ws.Range("M132").Value = -5753.249899999999
$ws.Range("N132").Value = -21759.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 179.6
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 212
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 636
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1106
$ws.Range("H68").Value = 193096.31
$ws.Range("I68").Value = 278342.62
$ws.Range("J68").Value = 1292.125
$ws.Range("K68").Value = 835027.86
$ws.Range("L68").Value = 3876.375
$ws.Range("M68").Value = -834216.86
$ws.Range("N68").Value = -5498.375
$ws.Range("H71").Value = 193096.31
$ws.Range("I71").Value = 278342.62
$ws.Range("J71").Value = 1292.125
$ws.Range("K71").Value = 2505083.58
$ws.Range("L71").Value = 11629.125
$ws.Range("M71").Value = -2501027.58
$ws.Range("N71").Value = -19741.125
$ws.Range("H107").Value = 734.35
$ws.Range("I107").Value = 683
$ws.Range("K107").Value = 2049
$ws.Range("M107").Value = -129
$ws.Range("H113").Value = 238789.42
$ws.Range("I113").Value = 500553.16
$ws.Range("J113").Value = 822.4091
$ws.Range("K113").Value = 1501659.48
$ws.Range("L113").Value = 2467.2273
$ws.Range("M113").Value = -1499489.48
$ws.Range("N113").Value = -6807.2273
$ws.Range("H131").Value = 2623.328
$ws.Range("I131").Value = 483.5
$ws.Range("J131").Value = 3519.0698
$ws.Range("K131").Value = 1450.5
$ws.Range("L131").Value = 10557.2094
$ws.Range("M131").Value = 3589.5
$ws.Range("N131").Value = -20637.2094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("H70").Value = 11551
$ws.Range("I70").Value = 16852
$ws.Range("J70").Value = 6250
$ws.Range("K70").Value = 16852
$ws.Range("L70").Value = 6250
$ws.Range("M70").Value = -16582
$ws.Range("N70").Value = -6790
$ws.Range("H73").Value = 11551
$ws.Range("I73").Value = 16852
$ws.Range("J73").Value = 6250
$ws.Range("K73").Value = 16852
$ws.Range("L73").Value = 6250
$ws.Range("M73").Value = -15916
$ws.Range("N73").Value = -8122
$ws.Range("H102").Value = 3649.9375
$ws.Range("I102").Value = 3490.818
$ws.Range("K102").Value = 3490.818
$ws.Range("M102").Value = -1868.818
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4514.143
$ws.Range("I7").Value = 4099.75
$ws.Range("K7").Value = 4099.75
$ws.Range("M7").Value = -3987.75
$ws.Range("H61").Value = 5833.3335
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -3298
$ws.Range("N61").Value = -7404
$ws.Range("H68").Value = 2285.3572
$ws.Range("I68").Value = 1921.6666
$ws.Range("J68").Value = 2940
$ws.Range("K68").Value = 1921.6666
$ws.Range("L68").Value = 2940
$ws.Range("M68").Value = -1172.6666
$ws.Range("N68").Value = -4438
$ws.Range("H71").Value = 2285.3572
$ws.Range("I71").Value = 1921.6666
$ws.Range("J71").Value = 2940
$ws.Range("K71").Value = 9608.333000000001
$ws.Range("L71").Value = 14700
$ws.Range("M71").Value = -5864.333000000001
$ws.Range("N71").Value = -22188
$ws.Range("H113").Value = 5833.3335
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -1330
$ws.Range("N113").Value = -11340
$ws.Range("H122").Value = 13338587
$ws.Range("I122").Value = 4864
$ws.Range("J122").Value = 20005448
$ws.Range("K122").Value = 14592
$ws.Range("L122").Value = 60016344
$ws.Range("M122").Value = -12142
$ws.Range("N122").Value = -60021244
$ws.Range("H126").Value = 4514.143
$ws.Range("I126").Value = 4099.75
$ws.Range("K126").Value = 12299.25
$ws.Range("M126").Value = -9829.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1289.2727
$ws.Range("I113").Value = 1050.909
$ws.Range("J113").Value = 1527.6364
$ws.Range("K113").Value = 3152.727
$ws.Range("L113").Value = 4582.9092
$ws.Range("M113").Value = -982.7270000000003
$ws.Range("N113").Value = -8922.9092
